# Apply the updated crypto price/volume figures (and the Chainlink/Cosmos row swap)
# captured by the commit. Values are prefixed with a leading apostrophe so Excel
# stores them as literal text (matching the original inlineStr cells) instead of
# auto-converting numeric-looking strings (e.g. "4.863", "0.9966") into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''30.065.00'
$ws.Range("E2").Value = '''  +0.22%  '
$ws.Range("D3").Value = '''1.882.03'
$ws.Range("E3").Value = '''  +0.69%  '
$ws.Range("D4").Value = '''0.9966'
$ws.Range("E4").Value = '''  -0.39%  '
$ws.Range("D5").Value = '''244.03'
$ws.Range("E5").Value = '''  -2.14%  '
$ws.Range("D6").Value = '''0.9966'
$ws.Range("E6").Value = '''  -0.45%  '
$ws.Range("D7").Value = '''0.4943'
$ws.Range("E7").Value = '''  -0.67%  '
$ws.Range("D8").Value = '''44.29'
$ws.Range("E8").Value = '''  -2.28%  '
$ws.Range("D9").Value = '''0.2920'
$ws.Range("E9").Value = '''  +3.04%  '
$ws.Range("D10").Value = '''0.06628'
$ws.Range("E10").Value = '''  +1.41%  '
$ws.Range("D11").Value = '''1.878.97'
$ws.Range("E11").Value = '''  +0.60%  '
$ws.Range("D12").Value = '''17.03'
$ws.Range("E12").Value = '''  +0.52%  '
$ws.Range("D13").Value = '''0.07199'
$ws.Range("E13").Value = '''  -0.25%  '
$ws.Range("D14").Value = '''0.6651'
$ws.Range("E14").Value = '''  +0.86%  '
$ws.Range("D15").Value = '''85.59'
$ws.Range("E15").Value = '''  +0.99%  '
$ws.Range("D16").Value = '''4.863'
$ws.Range("E16").Value = '''  +1.43%  '
$ws.Range("D17").Value = '''30.042.18'
$ws.Range("E17").Value = '''  +0.26%  '
$ws.Range("D18").Value = '''0.000007878'
$ws.Range("E18").Value = '''  +5.35%  '
$ws.Range("D19").Value = '''0.9974'
$ws.Range("E19").Value = '''  -0.43%  '
$ws.Range("D20").Value = '''12.81'
$ws.Range("E20").Value = '''  -0.01%  '
$ws.Range("D21").Value = '''2.121.16'
$ws.Range("E21").Value = '''  +0.44%  '
$ws.Range("D22").Value = '''0.9959'
$ws.Range("E22").Value = '''  -0.38%  '
$ws.Range("D23").Value = '''4.766'
$ws.Range("E23").Value = '''  +0.74%  '
$ws.Range("B24").Value = '''Chainlink'
$ws.Range("C24").Value = '''https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D24").Value = '''5.599'
$ws.Range("E24").Value = '''  +2.01%  '
$ws.Range("B25").Value = '''Cosmos'
$ws.Range("C25").Value = '''https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").Value = '''9.151'
$ws.Range("E25").Value = '''  +1.32%  '
$ws.Range("D26").Value = '''149.21'
$ws.Range("D27").Value = '''135.83'
$ws.Range("E27").Value = '''  +0.55%  '
$ws.Range("D28").Value = '''16.79'
$ws.Range("E28").Value = '''  +0.35%  '
$ws.Range("D29").Value = '''1.911'
$ws.Range("E29").Value = '''  -1.45%  '
$ws.Range("E30").Value = '''  -0.83%  '
$ws.Range("D31").Value = '''4.201'
$ws.Range("E31").Value = '''  -0.75%  '
$ws.Range("D32").Value = '''0.08649'
$ws.Range("E32").Value = '''  +0.49%  '
$ws.Range("D33").Value = '''3.959'
$ws.Range("E33").Value = '''  +2.11%  '
$ws.Range("D34").Value = '''0.04984'
$ws.Range("E34").Value = '''  -1.65%  '
$ws.Range("D35").Value = '''1.112'
$ws.Range("E35").Value = '''  -1.26%  '
$ws.Range("D36").Value = '''0.7022'
$ws.Range("E36").Value = '''  +2.85%  '
$ws.Range("D37").Value = '''2.653'
$ws.Range("E37").Value = '''  -1.48%  '
$ws.Range("D38").Value = '''2.212'
$ws.Range("E38").Value = '''  -5.28%  '
$ws.Range("D39").Value = '''2.695'
$ws.Range("E39").Value = '''  -1.33%  '
$ws.Range("D40").Value = '''0.9324'
$ws.Range("E40").Value = '''  -2.96%  '
$ws.Range("D41").Value = '''0.01643'
$ws.Range("E41").Value = '''  +1.01%  '
$ws.Range("D42").Value = '''5.991'
$ws.Range("E42").Value = '''  -1.99%  '
$ws.Range("D43").Value = '''0.9975'
$ws.Range("E43").Value = '''  -0.39%  '
$ws.Range("D44").Value = '''0.4207'
$ws.Range("E44").Value = '''  +0.64%  '
$ws.Range("D45").Value = '''101.87'
$ws.Range("E45").Value = '''  -1.75%  '
$ws.Range("D46").Value = '''7.582'
$ws.Range("E46").Value = '''  +2.13%  '
$ws.Range("D47").Value = '''0.1262'
$ws.Range("E47").Value = '''  +0.71%  '
$ws.Range("D48").Value = '''0.05711'
$ws.Range("E48").Value = '''  +1.63%  '
$ws.Range("D49").Value = '''32.50'
$ws.Range("E49").Value = '''  +0.28%  '
$ws.Range("D50").Value = '''8.254'
$ws.Range("E50").Value = '''  +0.25%  '
$ws.Range("D51").Value = '''0.3714'
$ws.Range("E51").Value = '''  -0.25%  '
